# Add the new "Sheet3" worksheet after "data" and populate it with the
# payment-distribution pivot data (year-month / day / fee_unit / payer).
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws3 = $wb.Worksheets.Add($null, $ws1)
$ws3.Name = "Sheet3"

# Header row
$ws3.Cells.Item(1, 1).Value = "year-month"
$ws3.Cells.Item(1, 2).Value = "day"
$ws3.Cells.Item(1, 3).Value = "fee_unit"
$ws3.Cells.Item(1, 4).Value = "payer"

# Data rows: "year-month|day|fee_unit|payer"
$rows = @(
  "2023-7|5|1|DONGLI LIU",
  "2023-7|4|2|EDWINALBERT IMHOFF",
  "2023-7|1|1|XIAOPENG FENG",
  "2023-6|30|1|TRACY CARRIER",
  "2023-6|6|1|DONGLI LIU",
  "2023-5|30|1|TRACY CARRIER",
  "2023-5|1|2|EDWINALBERT IMHOFF",
  "2023-5|1|1|DONGLI LIU",
  "2023-5|1|1|XIAOPENG FENG",
  "2023-5|1|1|SIDDHARTH SINHA",
  "2023-5|1|1|SIDDHARTH SINHA",
  "2023-5|1|1|SIDDHARTH SINHA",
  "2023-4|28|1|TRACY CARRIER",
  "2023-4|1|1|XIAOPENG FENG",
  "2023-3|31|1|TRACY CARRIER",
  "2023-3|31|2|EDWINALBERT IMHOFF",
  "2023-3|30|1|DONGLI LIU",
  "2023-3|2|1|DONGLI LIU",
  "2023-3|2|1|EDWINALBERT IMHOFF",
  "2023-3|2|1|EDWINALBERT IMHOFF",
  "2023-3|1|1|XIAOPENG FENG",
  "2023-3|1|1|TRACY CARRIER",
  "2023-2|24|1|SIDDHARTH SINHA",
  "2023-2|1|1|XIAOPENG FENG",
  "2023-2|1|2|EDWINALBERT IMHOFF",
  "2023-2|1|1|DONGLI LIU",
  "2023-1|31|1|TRACY CARRIER",
  "2023-1|2|1|XIAOPENG FENG",
  "2023-1|1|1|DONGLI LIU",
  "2022-12|31|1|TRACY CARRIER",
  "2022-12|27|1|SIDDHARTH SINHA",
  "2022-12|27|1|SIDDHARTH SINHA",
  "2022-12|5|1|DONGLI LIU",
  "2022-12|5|2|EDWINALBERT IMHOFF",
  "2022-12|1|1|XIAOPENG FENG",
  "2022-12|1|1|TRACY CARRIER",
  "2022-11|14|1|SIDDHARTH SINHA",
  "2022-11|1|1|XIAOPENG FENG",
  "2022-11|1|2|EDWINALBERT IMHOFF",
  "2022-11|1|1|DONGLI LIU",
  "2022-10|28|1|TRACY CARRIER",
  "2022-10|4|1|SIDDHARTH SINHA",
  "2022-10|2|1|XIAOPENG FENG",
  "2022-10|2|1|DONGLI LIU",
  "2022-9|30|1|TRACY CARRIER",
  "2022-9|29|4|EDWINALBERT IMHOFF",
  "2022-9|9|1|SIDDHARTH SINHA",
  "2022-9|1|1|XIAOPENG FENG",
  "2022-9|1|1|TRACY CARRIER",
  "2022-8|31|1|DONGLI LIU",
  "2022-8|1|1|SIDDHARTH SINHA",
  "2022-8|1|1|XIAOPENG FENG",
  "2022-7|31|2|EDWINALBERT IMHOFF",
  "2022-7|30|1|TRACY CARRIER",
  "2022-7|29|1|DONGLI LIU",
  "2022-7|1|1|XIAOPENG FENG",
  "2022-7|1|1|TRACY CARRIER",
  "2022-6|30|1|SIDDHARTH SINHA",
  "2022-6|29|2|EDWINALBERT IMHOFF",
  "2022-6|29|1|DONGLI LIU",
  "2022-6|3|1|SIDDHARTH SINHA",
  "2022-6|3|2|EDWINALBERT IMHOFF",
  "2022-6|2|1|XIAOPENG FENG",
  "2022-5|31|1|TRACY CARRIER",
  "2022-5|27|1|DONGLI LIU"
)

$r = 2
foreach ($line in $rows) {
  $parts = $line.Split("|")
  $ws3.Cells.Item($r, 1).Value = $parts[0]
  $ws3.Cells.Item($r, 2).Value = [int]$parts[1]
  $ws3.Cells.Item($r, 3).Value = [int]$parts[2]
  $ws3.Cells.Item($r, 4).Value = $parts[3]
  $r = $r + 1
}

# The "data" sheet is no longer the active tab; restore its lingering
# selection to the cell the author left it on.
$ws1.Activate()
$ws1.Range("A46").Select()

# Sheet3 becomes the active/selected sheet, frozen at the header row, with
# the view left on the cell the author was last looking at. Activate it
# LAST so it ends up as the workbook's active tab.
$ws3.Activate()
$ws3.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws3.Range("H12").Select()
$ws3.PageSetup.Orientation = 1

# Workbook-level defined names Excel wrote out when the payment-distribution
# chart was built from (and then trimmed down on) Sheet3's J:L columns.
$n = $wb.Names.Add('_xlchart.v5.0', '=Sheet3!$J$1')
$n.Visible = $false
$n = $wb.Names.Add('_xlchart.v5.1', '=Sheet3!$J$2:$J$66')
$n.Visible = $false
$n = $wb.Names.Add('_xlchart.v5.2', '=Sheet3!$K$1')
$n.Visible = $false
$n = $wb.Names.Add('_xlchart.v5.3', '=Sheet3!$K$2:$K$66')
$n.Visible = $false
$n = $wb.Names.Add('_xlchart.v5.4', '=Sheet3!$L$1')
$n.Visible = $false
$n = $wb.Names.Add('_xlchart.v5.5', '=Sheet3!$L$2:$L$66')
$n.Visible = $false
$n = $wb.Names.Add('_xlchart.v5.6', '=Sheet3!#REF!')
$n.Visible = $false
$n = $wb.Names.Add('_xlchart.v5.7', '=Sheet3!#REF!')
$n.Visible = $false
$n = $wb.Names.Add('_xlchart.v5.8', '=Sheet3!#REF!')
$n.Visible = $false
$n = $wb.Names.Add('_xlchart.v5.9', '=Sheet3!#REF!')
$n.Visible = $false
$n = $wb.Names.Add('_xlchart.v5.10', '=Sheet3!#REF!')
$n.Visible = $false
$n = $wb.Names.Add('_xlchart.v5.11', '=Sheet3!#REF!')
$n.Visible = $false
